# model2_df_results.xlsx regeneration:
# updates the R^2 / RMSE / U columns (C:E) for rows 2-10 with the new
# model run, and re-applies the per-column 3-colour gradient (green scale)
# that highlights RMSE (D) and U (E) the way the original report generator
# (pandas Styler -> openpyxl) did, now with font colour flipped to near-white
# on the darkest swatches so the numbers stay legible.

function ToComColor([string]$rgbHex) {
    $r = [Convert]::ToInt32($rgbHex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($rgbHex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($rgbHex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BLACK = ToComColor("000000")
$WHITE = ToComColor("F1F1F1")

# row => @{ C=...; D=...; E=...; DFill=...; EFill=...; DFont=...; EFont=... }
$rows = @{
  2  = @{ C = -9.8811;            D = 1.2113;            E = 2.9993; DFill="F7FCF5"; EFill="F7FCF5"; DFont=$BLACK; EFont=$BLACK }
  3  = @{ C = -4.4112;            D = 1.0743;            E = 2.3165; DFill="B1E0AB"; EFill="6ABF71"; DFont=$BLACK; EFont=$BLACK }
  4  = @{ C = -1.9079;            D = 0.9713000000000001; E = 2.0788; DFill="5AB769"; EFill="2D954D"; DFont=$BLACK; EFont=$BLACK }
  5  = @{ C = -0.3126;            D = 0.7853;             E = 1.6991; DFill="00441B"; EFill="00441B"; DFont=$WHITE; EFont=$WHITE }
  6  = @{ C = -0.3461;            D = 0.8802;             E = 1.942;  DFill="1C8540"; EFill="117B38"; DFont=$BLACK; EFont=$WHITE }
  7  = @{ C = -0.6104000000000001; D = 0.9487;            E = 2.3774; DFill="45AD5F"; EFill="7CC87C"; DFont=$BLACK; EFont=$BLACK }
  8  = @{ C = -0.7715;            D = 1.0027;             E = 2.4963; DFill="78C679"; EFill="9CD797"; DFont=$BLACK; EFont=$BLACK }
  9  = @{ C = -0.842;             D = 1.0302;             E = 2.5348; DFill="90D18D"; EFill="A7DBA0"; DFont=$BLACK; EFont=$BLACK }
  10 = @{ C = -0.9889;            D = 1.0784;             E = 2.6277; DFill="B5E1AE"; EFill="BCE4B5"; DFont=$BLACK; EFont=$BLACK }
}

foreach ($r in 2..10) {
    $row = $rows[$r]

    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Interior.Color = ToComColor($row.DFill)
    $dCell.Font.Color = $row.DFont

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Interior.Color = ToComColor($row.EFill)
    $eCell.Font.Color = $row.EFont
}
